$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): prefix with an apostrophe so Excel keeps the
# numeric-looking strings as text (matches the source inlineStr cells),
# then reset the style so the quote-prefix flag does not linger on the cell.
$ws.Range('D2').Value = "'67.703.23"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'2.493.20"
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').Value = "'586.94"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'176.65"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Value = "'0.516"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D14').Value = "'25.72"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = "'67.694.16"
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').Value = "'2.583.38"
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Value = "'10.98"
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = "'351.14"
$ws.Range('D20').Style = 'Normal'
$ws.Range('D26').Value = "'9.14"
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Value = "'2.622.45"
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Value = "'0.999"
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Value = "'0.0₃0903"
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Value = "'508.50"
$ws.Range('D30').Style = 'Normal'
$ws.Range('D35').Value = "'0.123"
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Value = "'162.85"
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Value = "'18.35"
$ws.Range('D38').Style = 'Normal'
$ws.Range('D41').Value = "'1.74"
$ws.Range('D41').Style = 'Normal'
$ws.Range('D43').Value = "'4.84"
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Value = "'2.40"
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Value = "'145.06"
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Value = "'3.51"
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Value = "'0.516"
$ws.Range('D47').Style = 'Normal'
$ws.Range('D50').Value = "'0.587"
$ws.Range('D50').Style = 'Normal'

# Volume(1h) column (E): plain text values (already non-numeric due to
# the %, leading "+"/"-", and padding spaces), so a direct assignment is safe.
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('E6').Value = '  +3.83%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.55%  '
$ws.Range('E9').Value = '  +3.90%  '
$ws.Range('E10').Value = '  +0.29%  '
$ws.Range('E11').Value = '  +2.94%  '
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('E17').Value = '  +3.14%  '
$ws.Range('E18').Value = '  +2.29%  '
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('E23').Value = '  +3.13%  '
$ws.Range('E24').Value = '  +1.81%  '
$ws.Range('E25').Value = '  -0.92%  '
$ws.Range('E26').Value = '  -1.67%  '
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('E29').Value = '  +0.73%  '
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('E31').Value = '  +1.75%  '
$ws.Range('E32').Value = '  +2.41%  '
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  +5.30%  '
$ws.Range('E36').Value = '  +2.60%  '
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('E39').Value = '  +0.61%  '
$ws.Range('E41').Value = '  +3.74%  '
$ws.Range('E43').Value = '  +1.27%  '
$ws.Range('E44').Value = '  +1.69%  '
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('E46').Value = '  +2.26%  '
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('E48').Value = '  +2.03%  '
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('E50').Value = '  +0.73%  '
$ws.Range('E51').Value = '  +0.56%  '
